$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new trade row (row 6) - repeater/all-night experiment trade result
$ws.Range("A6").Value = 42649.656180555554
$ws.Range("B6").Value = $false
$ws.Range("C6").Value = 9920.24
$ws.Range("D6").Value = 9949.09
$ws.Range("E6").Value = 104.82
$ws.Range("F6").Value = 105.43
$ws.Range("G6").Value = $true
$ws.Range("H6").Value = 0.57999999999999996
$ws.Range("I6").Value = $false

# Match the formatting of the existing date/boolean-short-sell columns (style index 1)
$ws.Range("A3").Copy()
$ws.Range("A6").PasteSpecial(-4122)

$ws.Range("G3").Copy()
$ws.Range("G6").PasteSpecial(-4122)

$excel.CutCopyMode = $false
